$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 906; this shifts the existing rows
# 906-923 down to 909-926 (and Excel copies the row-906 formatting,
# e.g. the date numeric format in column D, onto the new rows).
$ws.Rows("906:908").Insert()

# Row 906 - new weekly record (Extra quality)
$ws.Cells.Item(906, 1).Value  = 11
$ws.Cells.Item(906, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(906, 3).Value  = "Bíobío"
$ws.Cells.Item(906, 4).Value  = 45239
$ws.Cells.Item(906, 5).Value  = 8
$ws.Cells.Item(906, 6).Value  = 100112020
$ws.Cells.Item(906, 7).Value  = "Tomate"
$ws.Cells.Item(906, 8).Value  = "Larga vida"
$ws.Cells.Item(906, 9).Value  = "Extra"
$ws.Cells.Item(906, 10).Value = 50
$ws.Cells.Item(906, 11).Value = 18000
$ws.Cells.Item(906, 12).Value = 18000
$ws.Cells.Item(906, 13).Value = 18000
$ws.Cells.Item(906, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(906, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(906, 16).Value = 1000
$ws.Cells.Item(906, 17).Value = 18
$ws.Cells.Item(906, 18).Value = "Hortaliza"

# Row 907 - new weekly record (Primera quality)
$ws.Cells.Item(907, 1).Value  = 11
$ws.Cells.Item(907, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(907, 3).Value  = "Bíobío"
$ws.Cells.Item(907, 4).Value  = 45239
$ws.Cells.Item(907, 5).Value  = 8
$ws.Cells.Item(907, 6).Value  = 100112020
$ws.Cells.Item(907, 7).Value  = "Tomate"
$ws.Cells.Item(907, 8).Value  = "Larga vida"
$ws.Cells.Item(907, 9).Value  = "Primera"
$ws.Cells.Item(907, 10).Value = 200
$ws.Cells.Item(907, 11).Value = 16000
$ws.Cells.Item(907, 12).Value = 16000
$ws.Cells.Item(907, 13).Value = 16000
$ws.Cells.Item(907, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(907, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(907, 16).Value = 889
$ws.Cells.Item(907, 17).Value = 18
$ws.Cells.Item(907, 18).Value = "Hortaliza"

# Row 908 - new weekly record (Segunda quality)
$ws.Cells.Item(908, 1).Value  = 11
$ws.Cells.Item(908, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(908, 3).Value  = "Bíobío"
$ws.Cells.Item(908, 4).Value  = 45239
$ws.Cells.Item(908, 5).Value  = 8
$ws.Cells.Item(908, 6).Value  = 100112020
$ws.Cells.Item(908, 7).Value  = "Tomate"
$ws.Cells.Item(908, 8).Value  = "Larga vida"
$ws.Cells.Item(908, 9).Value  = "Segunda"
$ws.Cells.Item(908, 10).Value = 200
$ws.Cells.Item(908, 11).Value = 14000
$ws.Cells.Item(908, 12).Value = 14000
$ws.Cells.Item(908, 13).Value = 14000
$ws.Cells.Item(908, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(908, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(908, 16).Value = 778
$ws.Cells.Item(908, 17).Value = 18
$ws.Cells.Item(908, 18).Value = "Hortaliza"

Write-Host "done"
